# Complete test for Table_Drawdowns macro
# Adds a new test-case row (row 91) to the tests sheet:
#   C91 = "Table_Drawdowns_test"   (macro/test name)
#   B91 = "Test drawdown table"    (description)
#   A91 = "Table_Drawdowns"        (test identifier)
# and moves the active selection to the newly added row, matching the
# order in which the shared strings were appended to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(91, 3).Value = "Table_Drawdowns_test"
$ws.Cells.Item(91, 2).Value = "Test drawdown table"
$ws.Cells.Item(91, 1).Value = "Table_Drawdowns"

$ws.Range("A91").Select()
